$d = $word.ActiveDocument

$replacements = @(
    @{old = "86×54=4644"; new = "28×62=1736"},
    @{old = "12×41=492";  new = "65×37=2405"},
    @{old = "16×53=848";  new = "34×54=1836"},
    @{old = "25×23=575";  new = "81×44=3564"},
    @{old = "73×47=3431"; new = "43×22=946"},
    @{old = "86×71=6106"; new = "73×98=7154"},
    @{old = "15×60=900";  new = "22×57=1254"},
    @{old = "27×61=1647"; new = "39×81=3159"},
    @{old = "12×85=1020"; new = "27×12=324"},
    @{old = "82×51=4182"; new = "38×59=2242"},
    @{old = "53×15=795";  new = "14×49=686"},
    @{old = "71×55=3905"; new = "37×40=1480"},
    @{old = "69×58=4002"; new = "66×82=5412"},
    @{old = "54×16=864";  new = "53×81=4293"},
    @{old = "26×85=2210"; new = "74×51=3774"},
    @{old = "29×83=2407"; new = "72×25=1800"},
    @{old = "75×26=1950"; new = "25×72=1800"},
    @{old = "96×59=5664"; new = "50×20=1000"},
    @{old = "56×76=4256"; new = "74×69=5106"},
    @{old = "84×46=3864"; new = "59×36=2124"},
    @{old = "55×57=3135"; new = "39×14=546"},
    @{old = "69×69=4761"; new = "32×30=960"},
    @{old = "27×91=2457"; new = "91×39=3549"},
    @{old = "32×17=544";  new = "16×97=1552"},
    @{old = "40×20=800";  new = "63×86=5418"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
